# Apply daily-scrape update to the AIESEC Global Talent opportunities sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as TEXT (not auto-coerced to a number)
# without touching NumberFormat/Style (which would mint a brand-new cellXf).
# We build the literal via a throwaway formula cell, then copy/paste-special
# just the value - this mirrors how a user would "paste as text" in Excel
# and leaves the destination cell's style untouched.
$scratch = $ws.Cells.Item(1000, 1)
function Set-TextValue($cell, $text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

# New rows of data (rows 2-13). Old rows 14-16 are removed entirely.
$data = @(
    @("1330381", "https://aiesec.org/opportunity/global-talent/1330381", "Logistics Support Specialist", "Frankfurt am Main, Deutschland", "No", "12 applicants", "6 - 18 Months", "Greyfood GmbH"),
    @("1330185", "https://aiesec.org/opportunity/global-talent/1330185", "Front-End Web Developer", "Ciudad Juárez, Chihuahua, Mexico", "No", "0 applicants", "6 - 18 Months", "EP&O Corporation"),
    @("1329810", "https://aiesec.org/opportunity/global-talent/1329810", "[IMPACT FORTALEZA] Business Development - Marketplace Specialst", "Fortaleza - Zone 1, Fortaleza - Ceará, Brasil", "No", "18 applicants", "6 - 18 Months", "MAKRO MOVEIS E EQUIPAMENTOS MODULADOS LTDA"),
    @("1328965", "https://aiesec.org/opportunity/global-talent/1328965", "Account Manager (German level C1/C2 only)", "Assen, Nederland", "No", "15 applicants", "6 - 18 Months", "ICT Specialist"),
    @("1328328", "https://aiesec.org/opportunity/global-talent/1328328", "EMEA Senior Manager and Director Curriculum Coordinator", "40 Düsseldorf, Germany", "Yes", "130 applicants", "6 - 18 Months", "PwC Global Partnership"),
    @("1328113", "https://aiesec.org/opportunity/global-talent/1328113", "Junior Brand Manager - Trainee", "Bruxelles, Belgio", "No", "197 applicants", "6 - 18 Months", "UCB"),
    @("1327778", "https://aiesec.org/opportunity/global-talent/1327778", "Digital Content & Stakeholder Engagement Intern", "Colombo, Sri Lanka", "No", "18 applicants", "6 - 18 Months", "Solutions Ground (Pvt) Ltd"),
    @("1327300", "https://aiesec.org/opportunity/global-talent/1327300", "Language Specialist - French", "Colombo, Sri Lanka", "No", "28 applicants", "6 - 18 Months", "Aitken Spence Travels (Pvt) Ltd"),
    @("1327043", "https://aiesec.org/opportunity/global-talent/1327043", "Web Developer", "Sousse, Tunisia", "No", "22 applicants", "9 - 12 Weeks", "Progress Professional Center"),
    @("1325702", "https://aiesec.org/opportunity/global-talent/1325702", "Guest Relations Executive and Waitress", "Colombo, Sri Lanka", "No", "13 applicants", "3 - 6 Months", "Indian Kitchen PVT LTD"),
    @("1323468", "https://aiesec.org/opportunity/global-talent/1323468", "Sales Account Manager", "Cyberjaya, Selangor, Malaysia", "No", "49 applicants", "6 - 18 Months", "IX Telecom Sdn Bhd"),
    @("1307741", "https://aiesec.org/opportunity/global-talent/1307741", "Marketing Intern", "Cyberjaya, Selangor, Malaysia", "No", "231 applicants", "6 - 18 Months", "IX Telecom Sdn Bhd")
)

# Write data rows 2..13. Column A (opportunity id) must stay text even
# though it looks numeric, matching the source export.
$r = 2
foreach ($row in $data) {
    Set-TextValue $ws.Cells.Item($r, 1) $row[0]
    for ($c = 2; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

$scratch.ClearContents()

# Remove now-unused old rows 14-16.
$ws.Range("A14:H16").Delete()

# Highlight the "PREMIUM" flag of row 6 (now "Yes") with a yellow fill.
$ws.Range("E6").Interior.Color = 65535

# Adjust column widths for columns C, D, H per the new layout.
$ws.Columns.Item(3).ColumnWidth = 66
$ws.Columns.Item(4).ColumnWidth = 48
$ws.Columns.Item(8).ColumnWidth = 45
